# Finish first version of software specification
# - Switches the active tab to "FLASH log"
# - Fixes the uint16/17/18/19/20 typos in the FLASH log byte-size notes
#   (handled implicitly: we simply stop writing the old "Total bytes"
#   scratch table content, which drops those now-unused shared strings)
# - Replaces the old horizontal "day/month/.../relaysStatus" scratch table
#   on the "FLASH log" sheet with a proper vertical "Campo / Tamanho
#   [bytes] / Formato" table, including totals

$wb = $excel.ActiveWorkbook

$flash = $wb.Worksheets.Item("FLASH log")

# --- Clear out the old scratch table (row 3: C:P headers, row 4: C:P sizes,
#     row 5: C:P formats) but keep the "Total bytes" label (Q3) and the
#     running-total formula (Q4), matching the new blanked-out layout.
$flash.Range("C3:P3").ClearContents()
$flash.Range("B4").ClearContents()
$flash.Range("C4:P4").ClearContents()
$flash.Range("B5").ClearContents()
$flash.Range("C5:P5").ClearContents()

# --- New column widths for the vertical table ---
$flash.Columns.Item(2).ColumnWidth = 18.74   # column B ("Campo")
$flash.Columns.Item(3).ColumnWidth = 12.45   # column C ("Tamanho [bytes]")
for ($c = 4; $c -le 16; $c++) {
    $flash.Columns.Item($c).ColumnWidth = 11.31   # columns D:P ("Formato" + spares)
}

# --- Header row (row 9) ---
$headerRange = $flash.Range("B9:D9")
$flash.Rows.Item(9).RowHeight = 47.25

$flash.Range("B9").Value = "Campo"
$flash.Range("C9").Value = "Tamanho [bytes]"
$flash.Range("D9").Value = "Formato"

$headerRange.Font.Name = "Aptos Display"
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 12
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

# Tamanho [bytes] column header wraps, Campo/Formato headers don't.
$flash.Range("C9").WrapText = $true
$flash.Range("B9").WrapText = $false
$flash.Range("D9").WrapText = $false

# --- Data rows (10-23): Campo / Tamanho [bytes] / Formato ---
$fields = @(
    @{ Row = 10; Nome = "day";          Tamanho = 1; Formato = "uint8" },
    @{ Row = 11; Nome = "month";        Tamanho = 1; Formato = "uint8" },
    @{ Row = 12; Nome = "year";         Tamanho = 1; Formato = "uint8" },
    @{ Row = 13; Nome = "hour";         Tamanho = 1; Formato = "uint8" },
    @{ Row = 14; Nome = "min";          Tamanho = 1; Formato = "uint8" },
    @{ Row = 15; Nome = "sec";          Tamanho = 1; Formato = "uint8" },
    @{ Row = 16; Nome = "ambientTemp";  Tamanho = 1; Formato = "uint8" },
    @{ Row = 17; Nome = "ambientHum";   Tamanho = 1; Formato = "uint8" },
    @{ Row = 18; Nome = "soilHum1";     Tamanho = 2; Formato = "uint16" },
    @{ Row = 19; Nome = "soilHum2";     Tamanho = 2; Formato = "uint16" },
    @{ Row = 20; Nome = "soilHum3";     Tamanho = 2; Formato = "uint16" },
    @{ Row = 21; Nome = "soilHum4";     Tamanho = 2; Formato = "uint16" },
    @{ Row = 22; Nome = "soilHum5";     Tamanho = 2; Formato = "uint16" },
    @{ Row = 23; Nome = "relaysStatus"; Tamanho = 1; Formato = "uint8" }
)

foreach ($field in $fields) {
    $r = $field.Row
    $flash.Rows.Item($r).RowHeight = 21.75

    $flash.Cells.Item($r, 2).Value = $field.Nome
    $flash.Cells.Item($r, 3).Value = $field.Tamanho
    $flash.Cells.Item($r, 4).Value = $field.Formato

    $rowRange = $flash.Range($flash.Cells.Item($r, 2), $flash.Cells.Item($r, 4))
    $rowRange.Font.Name = "Lucida Console"
    $rowRange.Borders.LineStyle = 1
    $rowRange.HorizontalAlignment = -4108
    $rowRange.VerticalAlignment = -4108
}

# --- Totals rows ---
$flash.Rows.Item(24).RowHeight = 38.25
$flash.Range("B24").Value = "Total otimizado [bytes]"
$flash.Range("C24").Formula = "=SUM(C10:C23)"
$flash.Range("D24").Value2 = $null

$flash.Rows.Item(25).RowHeight = 31.5
$flash.Range("B25").Value = "Tamanho total ""alinhado"" [bytes]"
$flash.Range("C25").Formula = "=2*14"
$flash.Range("D25").Value2 = $null

$totalsLabelRange = $flash.Range("B24:B25")
$totalsLabelRange.Font.Name = "Aptos Display"
$totalsLabelRange.Font.Bold = $true
$totalsLabelRange.Font.Size = 12
$totalsLabelRange.Borders.LineStyle = 1
$totalsLabelRange.HorizontalAlignment = -4108
$totalsLabelRange.VerticalAlignment = -4108
$totalsLabelRange.WrapText = $true

$totalsValueRange = $flash.Range("C24:D25")
$totalsValueRange.Font.Name = "Lucida Console"
$totalsValueRange.Borders.LineStyle = 1
$totalsValueRange.HorizontalAlignment = -4108
$totalsValueRange.VerticalAlignment = -4108

# --- View state: select the new table and make "FLASH log" the active sheet/tab ---
$flash.Range("B9:D25").Select()
$flash.Activate()
$excel.ActiveWindow.ScrollRow = 7
